# Apply row shuffle for rows 137-151 (Jmt "Renkullmyren" fungal records)
# as described by the diff: several rows' taxon/coordinate data were
# permuted among each other while non-varying columns (location, dates,
# observer, etc.) stayed the same.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 137
$ws.Range("A137").Value = 112102104
$ws.Range("B137").Value = 90651
$ws.Range("D137").Value = "NT"
$ws.Range("E137").Value = 1968
$ws.Range("F137").Value = "Grantaggsvamp"
$ws.Range("G137").Value = "Bankera violascens"
$ws.Range("H137").Value = "(Alb. & Schwein. : Fr.) Pouzar"
$ws.Range("Q137").Value = 446883.2557033793
$ws.Range("R137").Value = 7032568.050871431

# Row 139
$ws.Range("A139").Value = 112102196
$ws.Range("B139").Value = 88956
$ws.Range("D139").Value = "VU"
$ws.Range("E139").Value = 5747
$ws.Range("F139").Value = "Läderdoftande fingersvamp"
$ws.Range("G139").Value = "Ramaria safraniolens"
$ws.Range("H139").Value = "Christian"
$ws.Range("Q139").Value = 446963.5003665166
$ws.Range("R139").Value = 7032564.927270759

# Row 140
$ws.Range("A140").Value = 112111606
$ws.Range("B140").Value = 88899
$ws.Range("D140").Value = "NT"
$ws.Range("E140").Value = 3286
$ws.Range("F140").Value = "Flattoppad klubbsvamp"
$ws.Range("G140").Value = "Clavariadelphus truncatus"
$ws.Range("H140").Value = "(Quél.) Donk"
$ws.Range("Q140").Value = 446863.3105291876
$ws.Range("R140").Value = 7032717.55581628

# Row 141
$ws.Range("A141").Value = 112104547
$ws.Range("B141").Value = 88946
$ws.Range("D141").Value = "VU"
$ws.Range("E141").Value = 256335
$ws.Range("F141").Value = "Taggfingersvamp"
$ws.Range("G141").Value = "Ramaria karstenii"
$ws.Range("H141").Value = "(Sacc. & P.Syd.) Corner"
$ws.Range("Q141").Value = 446688.4300003364
$ws.Range("R141").Value = 7032559.644946836

# Row 142
$ws.Range("A142").Value = 112103325
$ws.Range("B142").Value = 85313
$ws.Range("D142").Value = "NT"
$ws.Range("E142").Value = 3739
$ws.Range("F142").Value = "Persiljespindling"
$ws.Range("G142").Value = "Cortinarius sulfurinus"
$ws.Range("H142").Value = "Quél."
$ws.Range("Q142").Value = 446867.4804056262
$ws.Range("R142").Value = 7032725.549778301

# Row 143
$ws.Range("A143").Value = 112104270
$ws.Range("B143").Value = 85313
$ws.Range("D143").Value = "NT"
$ws.Range("E143").Value = 3739
$ws.Range("F143").Value = "Persiljespindling"
$ws.Range("G143").Value = "Cortinarius sulfurinus"
$ws.Range("H143").Value = "Quél."
$ws.Range("Q143").Value = 446732.1312419278
$ws.Range("R143").Value = 7032597.890938614

# Row 144
$ws.Range("A144").Value = 112102200
$ws.Range("B144").Value = 88899
$ws.Range("D144").Value = "NT"
$ws.Range("E144").Value = 3286
$ws.Range("F144").Value = "Flattoppad klubbsvamp"
$ws.Range("G144").Value = "Clavariadelphus truncatus"
$ws.Range("H144").Value = "(Quél.) Donk"
$ws.Range("Q144").Value = 446961.2729180742
$ws.Range("R144").Value = 7032565.860167116

# Row 145
$ws.Range("A145").Value = 112111486
$ws.Range("B145").Value = 82949
$ws.Range("D145").Value = "NT"
$ws.Range("E145").Value = 5589
$ws.Range("F145").Value = "Rödbrun klubbdyna"
$ws.Range("G145").Value = "Trichoderma nybergianum"
$ws.Range("H145").Value = "(T.Ulvinen & H.L.Chamb.) Jaklitsch & Voglmayr"
$ws.Range("Q145").Value = 446833.4140082744
$ws.Range("R145").Value = 7032727.011846939

# Row 146
$ws.Range("A146").Value = 112101773
$ws.Range("B146").Value = 86223
$ws.Range("D146").Value = "NT"
$ws.Range("E146").Value = 4412
$ws.Range("F146").Value = "Äggvaxskivling"
$ws.Range("G146").Value = "Hygrophorus karstenii"
$ws.Range("H146").Value = "Sacc. & Cub."
$ws.Range("Q146").Value = 446983.6614166541
$ws.Range("R146").Value = 7032942.217245953

# Row 147
$ws.Range("A147").Value = 112104266
$ws.Range("B147").Value = 88956
$ws.Range("D147").Value = "VU"
$ws.Range("E147").Value = 5747
$ws.Range("F147").Value = "Läderdoftande fingersvamp"
$ws.Range("G147").Value = "Ramaria safraniolens"
$ws.Range("H147").Value = "Christian"
$ws.Range("Q147").Value = 446732.1312419278
$ws.Range("R147").Value = 7032597.890938614

# Row 148
$ws.Range("A148").Value = 112102682
$ws.Range("B148").Value = 84820
$ws.Range("D148").Value = "VU"
$ws.Range("E148").Value = 275
$ws.Range("F148").Value = "Kejsarskivling"
$ws.Range("G148").Value = "Catathelasma imperiale"
$ws.Range("H148").Value = "(P.Karst.) Singer"
$ws.Range("Q148").Value = 447024.9334127782
$ws.Range("R148").Value = 7032671.870272635

# Row 149
$ws.Range("A149").Value = 112111498
$ws.Range("B149").Value = 88033
$ws.Range("D149").Value = "VU"
$ws.Range("E149").Value = 1599
$ws.Range("F149").Value = "Fjällfotad musseron"
$ws.Range("G149").Value = "Tricholoma olivaceotinctum"
$ws.Range("H149").Value = "Mort.Chr. & Heilm.-Claus."
$ws.Range("Q149").Value = 446860.139727794
$ws.Range("R149").Value = 7032742.69407742

# Row 150
$ws.Range("A150").Value = 112101944
$ws.Range("B150").Value = 88899
$ws.Range("D150").Value = "NT"
$ws.Range("E150").Value = 3286
$ws.Range("F150").Value = "Flattoppad klubbsvamp"
$ws.Range("G150").Value = "Clavariadelphus truncatus"
$ws.Range("H150").Value = "(Quél.) Donk"
$ws.Range("Q150").Value = 446857.8975496973
$ws.Range("R150").Value = 7032742.731334708

# Row 151
$ws.Range("A151").Value = 112104573
$ws.Range("B151").Value = 88033
$ws.Range("D151").Value = "VU"
$ws.Range("E151").Value = 1599
$ws.Range("F151").Value = "Fjällfotad musseron"
$ws.Range("G151").Value = "Tricholoma olivaceotinctum"
$ws.Range("H151").Value = "Mort.Chr. & Heilm.-Claus."
$ws.Range("Q151").Value = 446696.0171395433
$ws.Range("R151").Value = 7032530.399564721
